# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the symbol list
# with the latest scraped values. Cells are stored as text (prices/percents
# carry formatting such as "%"), so each value is written with a leading
# apostrophe to keep Excel from reinterpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.13"
$ws.Range("E2").Value = "'0.12%"
$ws.Range("D3").Value = "'27.01"
$ws.Range("E3").Value = "'-0.59%"
$ws.Range("D4").Value = "'4.718"
$ws.Range("E4").Value = "'-0.37%"
$ws.Range("D5").Value = "'0.06221"
$ws.Range("E5").Value = "'2.29%"
$ws.Range("D6").Value = "'6.750"
$ws.Range("E6").Value = "'1.23%"
$ws.Range("D7").Value = "'0.8510"
$ws.Range("E7").Value = "'0.39%"
$ws.Range("D8").Value = "'0.9144"
$ws.Range("E8").Value = "'-0.63%"
$ws.Range("D9").Value = "'0.1405"
$ws.Range("E9").Value = "'-0.13%"
$ws.Range("D10").Value = "'0.05075"
$ws.Range("E10").Value = "'2.70%"
$ws.Range("D11").Value = "'0.07081"
$ws.Range("E11").Value = "'-0.16%"
$ws.Range("E12").Value = "'-1.29%"
$ws.Range("D13").Value = "'0.09055"
$ws.Range("E13").Value = "'-0.30%"
$ws.Range("D14").Value = "'0.001528"
$ws.Range("E14").Value = "'-0.27%"
$ws.Range("D15").Value = "'0.0006177"
$ws.Range("E15").Value = "'1.87%"
$ws.Range("E16").Value = "'-1.87%"
$ws.Range("D18").Value = "'3.168"
$ws.Range("E18").Value = "'0.48%"
$ws.Range("D22").Value = "'4.081"
$ws.Range("E22").Value = "'-0.51%"
$ws.Range("E23").Value = "'0.05%"
$ws.Range("D24").Value = "'0.001201"
$ws.Range("E24").Value = "'-1.51%"
$ws.Range("D25").Value = "'0.004075"
$ws.Range("E25").Value = "'4.18%"
$ws.Range("E26").Value = "'0.02%"
$ws.Range("E27").Value = "'4.12%"
$ws.Range("D40").Value = "'0.03952"
$ws.Range("E40").Value = "'1.94%"
$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'0.04%"
$ws.Range("D42").Value = "'0.004135"
$ws.Range("E42").Value = "'0.12%"
$ws.Range("E43").Value = "'0.12%"
$ws.Range("E44").Value = "'-16.91%"
$ws.Range("E45").Value = "'-3.11%"
$ws.Range("E46").Value = "'0.02%"
$ws.Range("D48").Value = "'0.1995"
$ws.Range("E48").Value = "'47.42%"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("E50").Value = "'0.02%"
